$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 15.65
$ws.Range("E2").Value = 64.3
$ws.Range("F2").Value = 4.65
$ws.Range("N2").Value = 52.47848103381103

# Row 3
$ws.Range("D3").Value = 92479.06
$ws.Range("F3").Value = 1.79
$ws.Range("N3").Value = 52.47848103381103

# Row 4
$ws.Range("D4").Value = 278.97
$ws.Range("E4").Value = 48.2
$ws.Range("F4").Value = 5.28
$ws.Range("N4").Value = 52.47848103381103

# Row 5
$ws.Range("D5").Value = 12.47
$ws.Range("E5").Value = 47.6
$ws.Range("F5").Value = 12.2
$ws.Range("N5").Value = 52.47848103381103

# Row 6
$ws.Range("D6").Value = 188.96
$ws.Range("E6").Value = 40.6
$ws.Range("F6").Value = 7.58
$ws.Range("N6").Value = 52.47848103381103
